$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update enchantment description text in column B (rows 2-40)
$descriptions = @{
    2 = 'Increases mining speed while underwater'
    3 = 'Increases damage dealt to arthropods'
    4 = 'Makes the wearer unable to remove the armor'
    5 = 'Reduces damage taken from explosions'
    6 = 'Strikes lightning when thrown in a thunderstorm '
    7 = 'Increases swimming speed'
    8 = 'Increases mining speed'
    9 = 'Reduces fall damage taken'
    10 = 'Sets the target on fire'
    11 = 'Reduces damage taken from fire'
    12 = 'Sets shot arrows on fire'
    13 = 'Increases drop amount from some blocks'
    14 = 'Creates an ice path when walking on water'
    15 = 'Increases damage to aquatic mobs'
    16 = 'Prevents arrows from being consumed when shot'
    17 = 'Increases knockback dealt'
    18 = 'Increases amount of loot dropped by mobs'
    19 = 'Returns trident to thrower after it is thrown'
    20 = 'Increases chance to catch treasure'
    21 = 'Increases bite rate'
    22 = 'Repairs items in exchange for experience'
    23 = 'Shoots three projectiles at once'
    24 = 'Allows arrows to pierce through targets and shields'
    25 = 'Increases damage dealt with arrows'
    26 = 'Reduces damage taken from projectiles'
    27 = 'Reduces all types of damage taken'
    28 = 'Increases knockback dealt with arrows'
    29 = 'Reloads the crossbow quicker'
    30 = 'Extends breathing time underwater'
    31 = 'Launches the user when used in water'
    32 = 'Increases damage dealt to everything'
    33 = 'Makes drops block as they are'
    34 = 'Increases damage dealt to undead tagets'
    35 = 'Increases movement speed on soul blocks'
    36 = 'Increases sweep attack damage'
    37 = 'Increases speed while sneaking'
    38 = 'Deals damage to those who attack the wearer'
    39 = 'Increases durability'
    40 = 'Cursed item disappears upon death'
}

foreach ($row in $descriptions.Keys) {
    $ws.Range("B$row").Value = $descriptions[$row]
}

# Update the active cell selection to match the saved view state
$ws.Range("C41").Select()
